$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated stat values (sval data regenerated to filter save games)
$newB = 3.286832544864788
$newC = 1.655778082260271
$newD = 0.1494219747398047
$newE = 0.4942365360607697
$newG = 5.586269137925634

foreach ($row in 2..4) {
    $ws.Range("B$row").Value = $newB
    $ws.Range("C$row").Value = $newC
    $ws.Range("D$row").Value = $newD
    $ws.Range("E$row").Value = $newE
    $ws.Range("G$row").Value = $newG
}
